$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Network" label cells to "Network-maria"
$ws.Range("A6").Value = "Network-maria"
$ws.Range("A7").Value = "Network-maria"

# Widen column A to match the target OOXML width of ~22.42578125 characters
# (the engine quantizes column width to an internal pixel grid, so the
# supplied ColumnWidth is chosen to land on the closest achievable value)
$ws.Columns.Item(1).ColumnWidth = 21.666666667

# Update selection to C2
$ws.Range("C2").Select()
